$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PersonHealthInformation")
$ws.Range("A3:AN3").Copy($ws.Range("A4:AN4"))
$ws.Range("A4").Value = "testT4145_1"
$ws.Rows(4).Style = $ws.Rows(3).Style
